# Auto-generated Excel COM-interop script
# Applies updated market-board derived values (currentAveragePrice*, LevePrice*, LeveProfit*)
# to the Leve profit tables on each job sheet, per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2129.3333
$ws.Range("I15").Value = 2129.3333
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 6387.999899999999
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -6218.999899999999

$ws.Range("H19").Value = 649.7
$ws.Range("I19").Value = 700
$ws.Range("J19").Value = 637.125
$ws.Range("K19").Value = 700
$ws.Range("L19").Value = 637.125
$ws.Range("M19").Value = -525
$ws.Range("N19").Value = -987.125

$ws.Range("H28").Value = 1240.1538
$ws.Range("I28").Value = 528.125
$ws.Range("J28").Value = 2379.4
$ws.Range("K28").Value = 528.125
$ws.Range("L28").Value = 2379.4
$ws.Range("M28").Value = -43.125

$ws.Range("H41").Value = 136.125
$ws.Range("I41").Value = 87
$ws.Range("J41").Value = 152.5
$ws.Range("K41").Value = 87
$ws.Range("L41").Value = 152.5
$ws.Range("M41").Value = 353

$ws.Range("H137").Value = 3412.9
$ws.Range("I137").Value = 1867.2858
$ws.Range("J137").Value = 7019.3335
$ws.Range("K137").Value = 5601.857400000001
$ws.Range("L137").Value = 21058.0005
$ws.Range("M137").Value = -3051.857400000001
$ws.Range("N137").Value = -26158.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 3101.25
$ws.Range("I8").Value = 2202.5
$ws.Range("J8").Value = 4000
$ws.Range("K8").Value = 2202.5
$ws.Range("L8").Value = 4000
$ws.Range("M8").Value = -2058.5

$ws.Range("H12").Value = 1245
$ws.Range("I12").Value = 1245
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1245
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1072

$ws.Range("H14").Value = 1251
$ws.Range("I14").Value = 1251
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1251
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1076

$ws.Range("H16").Value = 10332.667
$ws.Range("I16").Value = 10332.667
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 10332.667
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -10045.667

$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = $null

$ws.Range("H74").Value = 2078
$ws.Range("I74").Value = 1613.6666
$ws.Range("J74").Value = 3471
$ws.Range("K74").Value = 1613.6666
$ws.Range("L74").Value = 3471
$ws.Range("M74").Value = -739.6666

$ws.Range("H77").Value = 2078
$ws.Range("I77").Value = 1613.6666
$ws.Range("J77").Value = 3471
$ws.Range("K77").Value = 8068.333000000001
$ws.Range("L77").Value = 17355
$ws.Range("M77").Value = -3700.333000000001

$ws.Range("H97").Value = 921.6
$ws.Range("I97").Value = 948.63635
$ws.Range("J97").Value = 723.3333
$ws.Range("K97").Value = 948.63635
$ws.Range("L97").Value = 723.3333
$ws.Range("M97").Value = -452.63635

$ws.Range("H132").Value = 6612.5884
$ws.Range("I132").Value = 3545.0908
$ws.Range("J132").Value = 12236.333
$ws.Range("K132").Value = 10635.2724
$ws.Range("L132").Value = 36708.999
$ws.Range("M132").Value = -8105.2724
$ws.Range("N132").Value = -41768.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 11819996
$ws.Range("I105").Value = 556972.7
$ws.Range("J105").Value = 62503600
$ws.Range("K105").Value = 556972.7
$ws.Range("L105").Value = 62503600
$ws.Range("M105").Value = -555225.7

$ws.Range("H134").Value = 2756.0625
$ws.Range("I134").Value = 1899.7142
$ws.Range("J134").Value = 3422.111
$ws.Range("K134").Value = 5699.142599999999
$ws.Range("L134").Value = 10266.333
$ws.Range("M134").Value = -3164.142599999999
$ws.Range("N134").Value = -15336.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3681053.8
$ws.Range("I31").Value = 2456.5454
$ws.Range("J31").Value = 5440383
$ws.Range("K31").Value = 2456.5454
$ws.Range("L31").Value = 5440383
$ws.Range("M31").Value = -2161.5454
$ws.Range("N31").Value = -5440973

$ws.Range("H34").Value = 3681053.8
$ws.Range("I34").Value = 2456.5454
$ws.Range("J34").Value = 5440383
$ws.Range("K34").Value = 2456.5454
$ws.Range("L34").Value = 5440383
$ws.Range("M34").Value = -2254.5454
$ws.Range("N34").Value = -5440787

$ws.Range("H58").Value = 2061.1853
$ws.Range("I58").Value = 1772.9
$ws.Range("J58").Value = 2884.8572
$ws.Range("K58").Value = 1772.9
$ws.Range("L58").Value = 2884.8572
$ws.Range("M58").Value = -1569.9
$ws.Range("N58").Value = -3290.8572

$ws.Range("H107").Value = 4167449.5
$ws.Range("I107").Value = 5000499.5
$ws.Range("J107").Value = 2198.5
$ws.Range("K107").Value = 5000499.5
$ws.Range("L107").Value = 2198.5
$ws.Range("M107").Value = -4998579.5
$ws.Range("N107").Value = -6038.5

$ws.Range("H132").Value = 14499428
$ws.Range("I132").Value = 6341.3
$ws.Range("J132").Value = 25647956
$ws.Range("K132").Value = 19023.9
$ws.Range("L132").Value = 76943868
$ws.Range("M132").Value = -16493.9
$ws.Range("N132").Value = -76948928

$ws.Range("H134").Value = 8747.833000000001
$ws.Range("I134").Value = 8747.833000000001
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 26243.499
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -23708.499

$ws.Range("H136").Value = 2061.1853
$ws.Range("I136").Value = 1772.9
$ws.Range("J136").Value = 2884.8572
$ws.Range("K136").Value = 5318.700000000001
$ws.Range("L136").Value = 8654.571599999999
$ws.Range("M136").Value = -2768.700000000001
$ws.Range("N136").Value = -13754.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 36001.4
$ws.Range("I74").Value = 40013
$ws.Range("J74").Value = 34998.5
$ws.Range("K74").Value = 120039
$ws.Range("L74").Value = 104995.5
$ws.Range("M74").Value = -118978
$ws.Range("N74").Value = -107117.5

$ws.Range("H77").Value = 36001.4
$ws.Range("I77").Value = 40013
$ws.Range("J77").Value = 34998.5
$ws.Range("K77").Value = 360117
$ws.Range("L77").Value = 314986.5
$ws.Range("M77").Value = -354813
$ws.Range("N77").Value = -325594.5

$ws.Range("H81").Value = 2573.75
$ws.Range("I81").Value = 3198
$ws.Range("J81").Value = 1949.5
$ws.Range("K81").Value = 9594
$ws.Range("L81").Value = 5848.5
$ws.Range("M81").Value = -8471
$ws.Range("N81").Value = -8094.5

$ws.Range("H84").Value = 2573.75
$ws.Range("I84").Value = 3198
$ws.Range("J84").Value = 1949.5
$ws.Range("K84").Value = 28782
$ws.Range("L84").Value = 17545.5
$ws.Range("M84").Value = -23166
$ws.Range("N84").Value = -28777.5

$ws.Range("H93").Value = 4437.5
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 4437.5
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 13312.5
$ws.Range("M93").Value = $null
$ws.Range("N93").Value = -17056.5

$ws.Range("H109").Value = 10124
$ws.Range("I109").Value = 17810
$ws.Range("J109").Value = 5000
$ws.Range("K109").Value = 53430
$ws.Range("L109").Value = 15000
$ws.Range("M109").Value = -52390
$ws.Range("N109").Value = -17080

$ws.Range("H122").Value = 1442.1578
$ws.Range("I122").Value = 512
$ws.Range("J122").Value = 1774.3572
$ws.Range("K122").Value = 4608
$ws.Range("L122").Value = 15969.2148
$ws.Range("M122").Value = -2158

$ws.Range("H137").Value = 18423.334
$ws.Range("I137").Value = 1513.3334
$ws.Range("J137").Value = 35333.332
$ws.Range("K137").Value = 4540.0002
$ws.Range("L137").Value = 105999.996
$ws.Range("M137").Value = 559.9997999999996
$ws.Range("N137").Value = -116199.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 5000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 5000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 5000
$ws.Range("M3").Value = $null
$ws.Range("N3").Value = -5224

$ws.Range("H7").Value = 2302.4614
$ws.Range("I7").Value = 1702.2858
$ws.Range("J7").Value = 3002.6667
$ws.Range("K7").Value = 1702.2858
$ws.Range("L7").Value = 3002.6667
$ws.Range("M7").Value = -1590.2858
$ws.Range("N7").Value = -3226.6667

$ws.Range("H13").Value = 3670.3333

$ws.Range("H15").Value = 5000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 5000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 5000
$ws.Range("M15").Value = $null
$ws.Range("N15").Value = -5340

$ws.Range("H16").Value = 1791.4667
$ws.Range("I16").Value = 694.4
$ws.Range("J16").Value = 3985.6
$ws.Range("K16").Value = 694.4
$ws.Range("L16").Value = 3985.6
$ws.Range("M16").Value = -524.4
$ws.Range("N16").Value = -4325.6

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = $null
$ws.Range("N17").Value = $null

$ws.Range("H100").Value = 3798.6875
$ws.Range("I100").Value = 2984.5
$ws.Range("J100").Value = 9498
$ws.Range("K100").Value = 2984.5
$ws.Range("L100").Value = 9498
$ws.Range("M100").Value = -2443.5

$ws.Range("H126").Value = 2302.4614
$ws.Range("I126").Value = 1702.2858
$ws.Range("J126").Value = 3002.6667
$ws.Range("K126").Value = 5106.857400000001
$ws.Range("L126").Value = 9008.000100000001
$ws.Range("M126").Value = -2636.857400000001
$ws.Range("N126").Value = -13948.0001

$ws.Range("H132").Value = 6285.5713
$ws.Range("I132").Value = 7666.3335
$ws.Range("J132").Value = 5250
$ws.Range("K132").Value = 22999.0005
$ws.Range("L132").Value = 15750
$ws.Range("M132").Value = -20469.0005
$ws.Range("N132").Value = -20810

$ws.Range("H136").Value = 2573
$ws.Range("I136").Value = 2573
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7719
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5169
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 3751
$ws.Range("I17").Value = 6499.5
$ws.Range("J17").Value = 1002.5
$ws.Range("K17").Value = 6499.5
$ws.Range("L17").Value = 1002.5
$ws.Range("M17").Value = -6327.5
$ws.Range("N17").Value = -1346.5

$ws.Range("H107").Value = 930.64703
$ws.Range("I107").Value = 930.64703
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2791.94109
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -871.9410899999998

$ws.Range("H132").Value = 2188.0344
$ws.Range("I132").Value = 1573
$ws.Range("J132").Value = 3802.5
$ws.Range("K132").Value = 4719
$ws.Range("L132").Value = 11407.5
$ws.Range("M132").Value = -2189

$ws.Range("H133").Value = 33999.5
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 33999.5
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 33999.5
$ws.Range("N133").Value = -44119.5

$ws.Range("H136").Value = 2343.7104
$ws.Range("I136").Value = 1759.4242
$ws.Range("J136").Value = 6200
$ws.Range("K136").Value = 5278.2726
$ws.Range("L136").Value = 18600
$ws.Range("M136").Value = -2728.2726
$ws.Range("N136").Value = -23700
